# Refresh the crypto price/volume snapshot (GitHub Actions scheduled update).
# Only cell values change; no rows/columns/styles are added or removed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.940.39"
$ws.Range("E2").Value = "  +0.98%  "
$ws.Range("D3").Value = "1.847.20"
$ws.Range("E3").Value = "  +1.06%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").Value = "'309.71"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.40%  "
$ws.Range("D6").Value = "'1.010"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("D7").Value = "'0.4778"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +2.55%  "
$ws.Range("D8").Value = "'0.3667"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.76%  "
$ws.Range("D9").Value = "'0.07224"
$ws.Range("D9").ClearFormats()
$ws.Range("D10").Value = "'0.9265"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.47%  "
$ws.Range("D11").Value = "'19.68"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.17%  "
$ws.Range("D12").Value = "'0.07710"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.77%  "
$ws.Range("D13").Value = "1.806.48"
$ws.Range("E13").Value = "  -1.21%  "
$ws.Range("E14").Value = "  +0.97%  "
$ws.Range("D15").Value = "'6.415"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.16%  "
$ws.Range("D16").Value = "'88.83"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.22%  "
$ws.Range("D17").Value = "'1.013"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("D18").Value = "'0.000008639"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("D20").Value = "26.969.73"
$ws.Range("E20").Value = "  +0.93%  "
$ws.Range("E21").Value = "  +2.38%  "
$ws.Range("D22").Value = "'5.057"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.93%  "
$ws.Range("E23").Value = "  +0.94%  "
$ws.Range("E24").Value = "  +0.48%  "
$ws.Range("D25").Value = "'152.37"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("D26").Value = "'18.18"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.40%  "
$ws.Range("D27").Value = "'1.994"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.81%  "
$ws.Range("D28").Value = "'114.12"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("D29").Value = "'4.945"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.34%  "
$ws.Range("D30").Value = "'0.08887"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.92%  "
$ws.Range("D31").Value = "'3.314"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +5.18%  "
$ws.Range("D32").Value = "'1.173"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.82%  "
$ws.Range("E33").Value = "  +1.30%  "
$ws.Range("D34").Value = "'4.492"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.15%  "
$ws.Range("D35").Value = "'2.720"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.58%  "
$ws.Range("D36").Value = "'1.122"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.91%  "
$ws.Range("E37").Value = "  +1.61%  "
$ws.Range("D38").Value = "'0.05267"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.60%  "
$ws.Range("D39").Value = "'2.984"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.86%  "
$ws.Range("D40").Value = "'0.5195"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.43%  "
$ws.Range("D41").Value = "'7.003"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.82%  "
$ws.Range("E42").Value = "  +0.67%  "
$ws.Range("E43").Value = "  +2.07%  "
$ws.Range("D44").Value = "'10.54"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +5.56%  "
$ws.Range("D45").Value = "'0.4730"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.27%  "
$ws.Range("E46").Value = "  +0.35%  "
$ws.Range("D47").Value = "'101.68"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.42%  "
$ws.Range("E48").Value = "  +2.48%  "
$ws.Range("D49").Value = "'66.02"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.12%  "
$ws.Range("D50").Value = "'0.06029"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.50%  "
$ws.Range("D51").Value = "'0.8864"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +4.27%  "

